$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("D4").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Interior.ThemeColor = 2
Write-Output "done"
